$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '62.009.30'
$ws.Range('E2').Value = '  +4.41%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.078.50'
$ws.Range('E3').Value = '  +2.55%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '580.24'
$ws.Range('E5').Value = '  +3.04%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '142.30'
$ws.Range('E6').Value = '  +2.93%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.069.12'
$ws.Range('E8').Value = '  +2.68%  '

$ws.Range('E9').Value = '  +1.38%  '

$ws.Range('E10').Value = '  +5.71%  '

$ws.Range('E11').Value = '  +10.99%  '

$ws.Range('E13').Value = '  +4.48%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '35.24'
$ws.Range('E14').Value = '  +4.45%  '

$ws.Range('E15').Value = '  +0.59%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.586.41'
$ws.Range('E16').Value = '  +2.57%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '7.27'
$ws.Range('E17').Value = '  +0.26%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.076.18'
$ws.Range('E18').Value = '  +2.62%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '61.930.28'
$ws.Range('E19').Value = '  +4.36%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '448.25'
$ws.Range('E20').Value = '  +4.15%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.05'
$ws.Range('E21').Value = '  +2.78%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.734'
$ws.Range('E22').Value = '  +1.96%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.44'
$ws.Range('E23').Value = '  +4.60%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.76'
$ws.Range('E24').Value = '  +3.13%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '81.70'
$ws.Range('E25').Value = '  +1.03%  '

$ws.Range('E26').Value = '  +0.23%  '

$ws.Range('E27').Value = '  +4.41%  '

$ws.Range('E28').Value = '  -0.02%  '

$ws.Range('E29').Value = '  +4.55%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.23'
$ws.Range('E30').Value = '  +5.62%  '

$ws.Range('E31').Value = '  +11.70%  '

$ws.Range('E32').Value = '  +14.11%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '26.84'
$ws.Range('E33').Value = '  +4.30%  '

$ws.Range('E34').Value = '  +4.79%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0₃0795'
$ws.Range('E35').Value = '  +4.01%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.05'
$ws.Range('E36').Value = '  +1.93%  '

$ws.Range('E37').Value = '  +5.89%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '50.20'
$ws.Range('E38').Value = '  +2.50%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.95'
$ws.Range('E39').Value = '  +8.34%  '

$ws.Range('E40').Value = '  +1.67%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '418.51'
$ws.Range('E41').Value = '  +4.28%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.908.82'
$ws.Range('E42').Value = '  +5.41%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0370'
$ws.Range('E43').Value = '  +5.11%  '

$ws.Range('E44').Value = '  +9.70%  '

$ws.Range('E45').Value = '  +0.39%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.14'
$ws.Range('E46').Value = '  +6.75%  '

$ws.Range('E47').Value = '  +0.02%  '

$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '124.08'
$ws.Range('E48').Value = '  +1.78%  '

$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '34.81'
$ws.Range('E49').Value = '  -2.21%  '

$ws.Range('E50').Value = '  +0.40%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '24.46'
$ws.Range('E51').Value = '  +4.32%  '
